# Updates cryptos list data (coin prices / hourly volume % changes) on Sheet1.
# Two coin-name swaps are included (rows 12/13: Polkadot <-> WrappedEther;
# rows 42/43: Quant <-> Maker) along with per-row price/volume refreshes.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("D2").Value = "26.110.01"
$ws.Range("E2").Value = "  -0.03%  "
# Row 3
$ws.Range("D3").Value = "1.652.02"
$ws.Range("E3").Value = "  -0.16%  "
# Row 4
$ws.Range("E4").Value = "  -0.44%  "
# Row 5
$ws.Range("D5").Value = "'218.38"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.57%  "
# Row 6
$ws.Range("D6").Value = "'0.5300"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.15%  "
# Row 7
$ws.Range("E7").Value = "  -0.30%  "
# Row 8
$ws.Range("D8").Value = "'0.2607"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.74%  "
# Row 9
$ws.Range("D9").Value = "'0.06316"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.78%  "
# Row 10
$ws.Range("D10").Value = "'20.39"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.45%  "
# Row 11
$ws.Range("D11").Value = "'0.07754"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.66%  "
# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.679.88"
$ws.Range("E12").Value = "  +0.54%  "
# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.469"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.51%  "
# Row 14
$ws.Range("D14").Value = "'0.5461"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.77%  "
# Row 15
$ws.Range("D15").Value = "0.0₅8134"
$ws.Range("E15").Value = "  -0.23%  "
# Row 16
$ws.Range("D16").Value = "'65.24"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.09%  "
# Row 17
$ws.Range("D17").Value = "26.114.94"
$ws.Range("E17").Value = "  -0.27%  "
# Row 18
$ws.Range("E18").Value = "  -0.31%  "
# Row 19
$ws.Range("D19").Value = "'4.546"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.97%  "
# Row 20
$ws.Range("D20").Value = "'193.38"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.95%  "
# Row 21
$ws.Range("D21").Value = "'10.04"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.45%  "
# Row 22
$ws.Range("D22").Value = "'5.976"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.06%  "
# Row 23
$ws.Range("D23").Value = "'1.003"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.49%  "
# Row 24
$ws.Range("D24").Value = "'139.96"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.02%  "
# Row 25
$ws.Range("E25").Value = "  +0.84%  "
# Row 26
$ws.Range("D26").Value = "'7.253"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.21%  "
# Row 27
$ws.Range("E27").Value = "  +0.36%  "
# Row 28
$ws.Range("E28").Value = "  +1.84%  "
# Row 29
$ws.Range("D29").Value = "'0.05909"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.73%  "
# Row 30
$ws.Range("E30").Value = "  +0.04%  "
# Row 31
$ws.Range("D31").Value = "'3.501"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.92%  "
# Row 32
$ws.Range("D32").Value = "'3.240"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.75%  "
# Row 33
$ws.Range("E33").Value = "  -5.06%  "
# Row 34
$ws.Range("D34").Value = "'2.411"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.05%  "
# Row 35
$ws.Range("D35").Value = "'0.9444"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.79%  "
# Row 36
$ws.Range("D36").Value = "'2.755"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.80%  "
# Row 37
$ws.Range("D37").Value = "'0.5645"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.74%  "
# Row 38
$ws.Range("E38").Value = "  +1.55%  "
# Row 39
$ws.Range("D39").Value = "'5.834"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.42%  "
# Row 40
$ws.Range("D40").Value = "'0.8434"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.95%  "
# Row 41
$ws.Range("E41").Value = "  -0.15%  "
# Row 42
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "1.008.20"
$ws.Range("E42").Value = "  -2.36%  "
# Row 43
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "'100.73"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.43%  "
# Row 44
$ws.Range("D44").Value = "1.797.29"
$ws.Range("E44").Value = "  -0.18%  "
# Row 45
$ws.Range("D45").Value = "'56.82"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.08%  "
# Row 46
$ws.Range("E46").Value = "  -2.26%  "
# Row 47
$ws.Range("D47").Value = "'1.008"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.63%  "
# Row 48
$ws.Range("E48").Value = "  +1.32%  "
# Row 49
$ws.Range("E49").Value = "  +1.62%  "
# Row 50
$ws.Range("D50").Value = "'0.05149"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.60%  "
# Row 51
$ws.Range("D51").Value = "'7.801"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.91%  "
